$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cumulative-death rows (dates 2021-10-20 .. 2021-11-02), appended
# after the existing data which ends at row 369 (date serial 44488).
$newRows = @(
    @(370, 44489, 12883, 2510, 15393),
    @(371, 44490, 12886, 2511, 15397),
    @(372, 44491, 12895, 2511, 15406),
    @(373, 44492, 12903, 2515, 15418),
    @(374, 44493, 12917, 2515, 15432),
    @(375, 44494, 12935, 2516, 15451),
    @(376, 44495, 12957, 2518, 15475),
    @(377, 44496, 12977, 2521, 15498),
    @(378, 44497, 13000, 2521, 15521),
    @(379, 44498, 13018, 2525, 15543),
    @(380, 44499, 13034, 2529, 15563),
    @(381, 44500, 13045, 2529, 15574),
    @(382, 44501, 13076, 2529, 15605),
    @(383, 44502, 13112, 2533, 15645)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
}
